$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 135, shifting existing rows 135-216 down to 136-217
$ws.Rows(135).Insert()

# Populate the newly inserted row 135 with the new data record
$ws.Cells.Item(135, 1).Value = 11
$ws.Cells.Item(135, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(135, 3).Value = "Bíobío"
$ws.Cells.Item(135, 4).Value = 45097
$ws.Cells.Item(135, 5).Value = 8
$ws.Cells.Item(135, 6).Value = "Fruta"
$ws.Cells.Item(135, 7).Value = 100102
$ws.Cells.Item(135, 8).Value = "Cítricos"
$ws.Cells.Item(135, 9).Value = 100102004
$ws.Cells.Item(135, 10).Value = "Mandarina"
$ws.Cells.Item(135, 11).Value = "Clementina"
$ws.Cells.Item(135, 12).Value = "Primera"
$ws.Cells.Item(135, 13).Value = 100
$ws.Cells.Item(135, 14).Value = 9000
$ws.Cells.Item(135, 15).Value = 10000
$ws.Cells.Item(135, 16).Value = 9500
$ws.Cells.Item(135, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(135, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(135, 19).Value = 528
$ws.Cells.Item(135, 20).Value = 18
